$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12-30 shift down to 13-31.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new data record.
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "Macroferia Regional de Talca"
$ws.Range("C12").Value = "Maule"
$ws.Range("D12").Value = 44725
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino dulce"
$ws.Range("H12").Value = "Cultivar IV Región"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = '$/bandeja 18 kilos'
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 833
$ws.Range("Q12").Value = 18
$ws.Range("R12").Value = "Hortaliza"
